$wb = $excel.ActiveWorkbook

# --- Shared literal strings (reused across sheets, matches existing shared-string pool) ---
$neo4jUrl   = 'Neo4j_URL:'
$neo4jUrlV  = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$userName   = 'User_name:'
$userNameV  = 'neo4j'
$pwdLbl     = 'PWD:'
$pwdV       = 'icdcDBneo4j0'
$cypherLbl  = 'Cypher:'
$cypherV    = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.sex IN [''Spayed female''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$outputLbl  = 'Output:'
$outputV    = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC06_Canine_Filter_Gender-SpayedFemale_Neo4jData.xlsx'
$statCypherV = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.sex IN [''Spayed female'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# --- Add new sheet: CypherOutput_Message (copy of Message content) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sCOM = $wb.Worksheets.Add($null, $lastSheet)
$sCOM.Name = "CypherOutput_Message"
$sCOM.Range("A1").Value = $neo4jUrl
$sCOM.Range("A2").Value = $neo4jUrlV
$sCOM.Range("A3").Value = $userName
$sCOM.Range("A4").Value = $userNameV
$sCOM.Range("A5").Value = $pwdLbl
$sCOM.Range("A6").Value = $pwdV
$sCOM.Range("A7").Value = $cypherLbl
$sCOM.Range("A8").Value = $cypherV
$sCOM.Range("A9").Value = $outputLbl
$sCOM.Range("A10").Value = $outputV

# --- Add new sheet: StatOutput (counts table) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sStat = $wb.Worksheets.Add($null, $lastSheet)
$sStat.Name = "StatOutput"
$sStat.Range("A1").Value = 'number_of_files'
$sStat.Range("B1").Value = 'number_of_sample'
$sStat.Range("C1").Value = 'number_of_cases'
$sStat.Range("D1").Value = 'number_of_study'
# Leading apostrophe forces these numeric-looking values to be stored as
# text (matching the source data, which keeps counts as text cells).
$sStat.Range("A2").Value = "'0"
$sStat.Range("B2").Value = "'0"
$sStat.Range("C2").Value = "'28"
$sStat.Range("D2").Value = "'1"

# --- Add new sheet: StatOutput_Message (Message content twice, second Cypher replaced with stats query) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sStatMsg = $wb.Worksheets.Add($null, $lastSheet)
$sStatMsg.Name = "StatOutput_Message"
$sStatMsg.Range("A1").Value = $neo4jUrl
$sStatMsg.Range("A2").Value = $neo4jUrlV
$sStatMsg.Range("A3").Value = $userName
$sStatMsg.Range("A4").Value = $userNameV
$sStatMsg.Range("A5").Value = $pwdLbl
$sStatMsg.Range("A6").Value = $pwdV
$sStatMsg.Range("A7").Value = $cypherLbl
$sStatMsg.Range("A8").Value = $cypherV
$sStatMsg.Range("A9").Value = $outputLbl
$sStatMsg.Range("A10").Value = $outputV
$sStatMsg.Range("A11").Value = $neo4jUrl
$sStatMsg.Range("A12").Value = $neo4jUrlV
$sStatMsg.Range("A13").Value = $userName
$sStatMsg.Range("A14").Value = $userNameV
$sStatMsg.Range("A15").Value = $pwdLbl
$sStatMsg.Range("A16").Value = $pwdV
$sStatMsg.Range("A17").Value = $cypherLbl
$sStatMsg.Range("A18").Value = $statCypherV
$sStatMsg.Range("A19").Value = $outputLbl
$sStatMsg.Range("A20").Value = $outputV

$sCOM.Range("A1").Select()
